{"js": "// Replace the multiplication-table answer strings throughout the document.\n// Each old value is unique in the document, so a simple search+replace per\n// pair is sufficient and keeps the original run formatting (rFonts/sz) intact\n// because we only rewrite the text inside the matched range.\nconst replacements = [\n  [\"89\u00d735=3115\", \"29\u00d719=551\"],\n  [\"97\u00d723=2231\", \"24\u00d743=1032\"],\n  [\"33\u00d751=1683\", \"75\u00d726=1950\"],\n  [\"91\u00d717=1547\", \"56\u00d775=4200\"],\n  [\"42\u00d745=1890\", \"35\u00d776=2660\"],\n  [\"13\u00d749=637\", \"35\u00d731=1085\"],\n  [\"85\u00d721=1785\", \"15\u00d781=1215\"],\n  [\"42\u00d712=504\", \"27\u00d771=1917\"],\n  [\"19\u00d738=722\", \"41\u00d781=3321\"],\n  [\"35\u00d717=595\", \"24\u00d746=1104\"],\n  [\"70\u00d788=6160\", \"52\u00d713=676\"],\n  [\"67\u00d743=2881\", \"44\u00d743=1892\"],\n  [\"76\u00d788=6688\", \"96\u00d788=8448\"],\n  [\"63\u00d752=3276\", \"34\u00d730=1020\"],\n  [\"65\u00d783=5395\", \"29\u00d766=1914\"],\n  [\"12\u00d786=1032\", \"95\u00d713=1235\"],\n  [\"38\u00d733=1254\", \"97\u00d748=4656\"],\n  [\"55\u00d716=880\", \"29\u00d798=2842\"],\n  [\"86\u00d783=7138\", \"94\u00d717=1598\"],\n  [\"36\u00d711=396\", \"53\u00d788=4664\"],\n  [\"81\u00d729=2349\", \"20\u00d745=900\"],\n  [\"18\u00d798=1764\", \"41\u00d763=2583\"],\n  [\"43\u00d752=2236\", \"68\u00d747=3196\"],\n  [\"93\u00d724=2232\", \"94\u00d796=9024\"],\n  [\"32\u00d764=2048\", \"65\u00d794=6110\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the multiplication-table answer strings throughout the document.\n# Each old value is unique in the document, so a simple Find/Replace per pair\n# is sufficient and leaves the surrounding run formatting (rFonts/sz) intact.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"89\u00d735=3115\", \"29\u00d719=551\"),\n    @(\"97\u00d723=2231\", \"24\u00d743=1032\"),\n    @(\"33\u00d751=1683\", \"75\u00d726=1950\"),\n    @(\"91\u00d717=1547\", \"56\u00d775=4200\"),\n    @(\"42\u00d745=1890\", \"35\u00d776=2660\"),\n    @(\"13\u00d749=637\", \"35\u00d731=1085\"),\n    @(\"85\u00d721=1785\", \"15\u00d781=1215\"),\n    @(\"42\u00d712=504\", \"27\u00d771=1917\"),\n    @(\"19\u00d738=722\", \"41\u00d781=3321\"),\n    @(\"35\u00d717=595\", \"24\u00d746=1104\"),\n    @(\"70\u00d788=6160\", \"52\u00d713=676\"),\n    @(\"67\u00d743=2881\", \"44\u00d743=1892\"),\n    @(\"76\u00d788=6688\", \"96\u00d788=8448\"),\n    @(\"63\u00d752=3276\", \"34\u00d730=1020\"),\n    @(\"65\u00d783=5395\", \"29\u00d766=1914\"),\n    @(\"12\u00d786=1032\", \"95\u00d713=1235\"),\n    @(\"38\u00d733=1254\", \"97\u00d748=4656\"),\n    @(\"55\u00d716=880\", \"29\u00d798=2842\"),\n    @(\"86\u00d783=7138\", \"94\u00d717=1598\"),\n    @(\"36\u00d711=396\", \"53\u00d788=4664\"),\n    @(\"81\u00d729=2349\", \"20\u00d745=900\"),\n    @(\"18\u00d798=1764\", \"41\u00d763=2583\"),\n    @(\"43\u00d752=2236\", \"68\u00d747=3196\"),\n    @(\"93\u00d724=2232\", \"94\u00d796=9024\"),\n    @(\"32\u00d764=2048\", \"65\u00d794=6110\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute([ref]$oldText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2) | Out-Null\n}\n"}
